$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Poroto granado" at Vega Modelo de
# Temuco. It belongs right above the existing row 24 (date-sorted data),
# so insert a new row there - this pushes the old rows 24-79 down to 25-80,
# which is exactly the shift seen throughout the diff.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new record's data.
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44624
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 100112030
$ws.Range("G24").Value = "Poroto granado"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 30
$ws.Range("K24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 25000
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Región de La Araucanía"
$ws.Range("P24").Value = 1000
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
